# The author replaced the full image URLs in column E (image_url) with
# short local filenames (tag1.jpg / tag2.jpg / tag3.jpg) for rows 2-4,
# while the underlying hyperlink targets for those cells are unchanged.
# Finally the active-cell selection moved to E5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "tag1.jpg"
$ws.Range("E3").Value = "tag2.jpg"
$ws.Range("E4").Value = "tag3.jpg"

$ws.Range("E5").Select()
